$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.451.89"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "3.106.88"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "388.62"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").Value = "103.88"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").Value = "37.17"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "0.0860"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "3.602.27"
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "18.57"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "7.85"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "3.115.76"
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("D17").Value = "1.00"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("D18").Value = "10.62"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").Value = "51.531.39"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "3.26"
$ws.Range("E20").Value = "  +5.91%  "
$ws.Range("D21").Value = "12.52"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "0.0₃0968"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").Value = "70.30"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "267.02"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "3.18"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "8.06"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("D27").Value = "27.36"
$ws.Range("E27").Value = "  +4.05%  "
$ws.Range("D28").Value = "7.19"
$ws.Range("E28").Value = "  -4.88%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -5.08%  "
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("D32").Value = "10.43"
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("D33").Value = "35.89"
$ws.Range("E33").Value = "  +4.89%  "
$ws.Range("D34").Value = "0.0478"
$ws.Range("E34").Value = "  +6.56%  "
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").Value = "50.04"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "3.39"
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").Value = "0.290"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "129.47"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.87"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "16.60"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").Value = "3.81"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").Value = "22.15"
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D47").Value = "2.49"
$ws.Range("E47").Value = "  +4.28%  "
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("D49").Value = "2.081.06"
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("D50").Value = "0.0333"
$ws.Range("E50").Value = "  +4.00%  "
$ws.Range("D51").Value = "0.923"
$ws.Range("E51").Value = "  +18.06%  "
